# Regenerate orders with updated distance/size codes.
#
# The stimulus-condition labels encoded in several columns (Condition,
# Filename_Left, Filename_Right, Distance, Size) use tokens like "D51",
# "D64", "D80" (viewing distance) and "S30" (stimulus size). This edit
# renumbers those tokens:
#   D51 -> D55
#   D64 -> D69
#   D80 -> D86
#   S30 -> S31
#
# The substitution is applied uniformly to every textual cell value in
# the used range, which naturally covers every column where the tokens
# appear (Condition, Filename_Left, Filename_Right, Distance, Size) and
# leaves numeric/boolean columns (Trial, Duration_Seconds, Is_Repeat,
# Block, ConditionID) and unrelated text (Face, NULL*, headers) unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$firstCol = $used.Column
$lastRow = $firstRow + $used.Rows.Count - 1
$lastCol = $firstCol + $used.Columns.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -is [string]) {
            $nv = $v -replace "D51", "D55" -replace "D64", "D69" -replace "D80", "D86" -replace "S30", "S31"
            if ($nv -ne $v) {
                $cell.Value2 = $nv
            }
        }
    }
}
